# Append a new numbered list item ("大后天开心到死") right after the
# existing "后天开心起飞" list item, matching its paragraph/list formatting
# (style "a3", numPr ilvl 0 / numId 1, ind firstLineChars 0).

$d = $word.ActiveDocument

# Locate the paragraph that ends with "后天开心起飞" so the new item is
# inserted in the right spot even if the document changes slightly.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*后天开心起飞*") {
        $target = $p
    }
}

if ($target -eq $null) {
    throw "Could not find paragraph containing '后天开心起飞'"
}

# Inserting a paragraph after an existing list item copies its paragraph
# formatting (style + numbering) onto the new paragraph automatically.
$newRange = $target.Range.InsertParagraphAfter()

# The newly created paragraph is now the one immediately following $target.
$newPara = $target.Next()
$newPara.Range.Text = "大后天开心到死"
